$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 9, 10, 11, 12, 13, 17, 18, 19, 20, 21, 22)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "nan"
}
